# Heineman's code package: the "toBeDeleted" entity row (row 6, which held
# the now-removed "toBeDeleted" label) is deleted outright. Excel shifts
# every row below it up by one, and the orphaned "toBeDeleted" shared
# string is dropped from the workbook on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").EntireRow.Delete()

# Leave the selection where the author last left it.
$ws.Range("D3").Select()
